$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.479.21"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.939.08"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.26"
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.59"
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.736"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").Value = "  +4.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.12"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.75"
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.569.97"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.72"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.956.50"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.04"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.574.26"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.11"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.49"
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.91"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.87"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.46"
$ws.Range("E25").Value = "  +13.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.18"
$ws.Range("E26").Value = "  +16.36%  "
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.04"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.84"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "716.56"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.70"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("E34").Value = "  +14.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.07"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.15"
$ws.Range("E36").Value = "  +13.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "61.12"
$ws.Range("E37").Value = "  +5.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.397"
$ws.Range("E39").Value = "  +18.04%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  +14.80%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0484"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("E44").Value = "  +6.40%  "
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.16"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0347"
$ws.Range("E49").Value = "  +34.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.13"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.16"
$ws.Range("E51").Value = "  -0.02%  "
